$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 125 (pushes the old row 125..394 down to 126..395,
# and extends the used range to A1:R395).
$ws.Rows.Item(125).Insert()

# Populate the newly inserted row 125 with its values.
$ws.Range("A125").Value2 = 3
$ws.Range("B125").Value2 = "Femacal de La Calera"
$ws.Range("C125").Value2 = "Coquimbo"
$ws.Range("D125").Value2 = 44935
$ws.Range("E125").Value2 = 5
$ws.Range("F125").Value2 = 100112039
$ws.Range("G125").Value2 = "Ciboulette"
$ws.Range("H125").Value2 = "Sin especificar"
$ws.Range("I125").Value2 = "Primera"
$ws.Range("J125").Value2 = 50
$ws.Range("K125").Value2 = 2000
$ws.Range("L125").Value2 = 2000
$ws.Range("M125").Value2 = 2000
$ws.Range("N125").Value2 = "`$/docena de atados"
$ws.Range("O125").Value2 = "Provincia de Quillota"
$ws.Range("P125").Value2 = 667
$ws.Range("Q125").Value2 = 3
$ws.Range("R125").Value2 = "Hortaliza"

# Make sure the date cell keeps the date number format used by the rest of
# column D (the Insert() above already propagates the neighbouring style,
# but set it explicitly to be safe).
$ws.Range("D125").NumberFormat = "YYYY-MM-DD HH:MM:SS"
